# Commit: "Added a zip file + updated rules/data1/data2"
#
# The microplastic image filenames in column Q no longer live inside the
# "MicroplasticImages.zip" archive path - the zip prefix is stripped so the
# sheet stores bare filenames.
#
# Also update the sheet's remembered cell selection (cosmetic UI state that
# Excel persists with the worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = "ALGALITA_CW_3_above500_30.jpeg"
$ws.Range("Q3").Value = "B_DW_3_above500_96.jpeg"
$ws.Range("Q4").Value = "CC_CW_1_20-250_176.jpeg"

$ws.Range("G9").Select()
